$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (pushes old rows 13-21 down to 14-22).
$ws.Rows.Item(13).Insert()

# The inserted row copied column A's bold style into A13; since the target
# layout has no entry at all in A13, reset its style and clear it.
$ws.Range("A13").Style = "Normal"
$ws.Range("A13").ClearContents()

# Give B13/C13 the same look as the row below (normal / red wrap-text column
# styles) by copying formats, then fill in the teacher's name.
$ws.Range("B14").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$teacherText = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("B13").Value = $teacherText
$ws.Range("C13").Value = $teacherText

# Row 10 (Objetivos:): replace the placeholder text with the real course
# objectives description.
$objetivosText = "Apresentar o conjunto de métodos de investigação de superfície e de sub-superfície (diretos e indiretos) utilizados nas caracterizações geológico-geotécnicas que envolvem o meio ambiente. Estabelecer análise crítica que possibilite a escolha e a utilização adequadas das técnicas de investigações disponíveis visando o estudo dos diversos tipos de problemas ambientais."
$ws.Range("B10").Value = $objetivosText
$ws.Range("C10").Value = $objetivosText

# Row 14 (Programa resumido:): replace "Semestral" with the real short
# syllabus summary.
$resumidoText = "Investigação de superfície e sub-superfícies; Técnicas e equipamentos mais adequados."
$ws.Range("B14").Value = $resumidoText
$ws.Range("C14").Value = $resumidoText

# Row 16 (Programa:): replace the erroneous date text with the full
# syllabus.
$programaText = "Introdução, conceitos e objetivos; Seqüência de estudos rotineiros (usuais); Estudos Corretivos e Preventivos; Investigação de Superfície; Investigação de Sub-superfície - Aplicações / Limitações - Métodos Diretos; Métodos Indiretos (Geofísicos); Métodos Sísmicos; Métodos Elétricos e Eletromagnético; Ensaios em Furos de Sondagem; Ensaios com Traçadores; Instrumentação Hidráulica e Mecânica. Estudo de caso."
$ws.Range("B16").Value = $programaText
$ws.Range("C16").Value = $programaText

# Row 19 (Método:): replace the erroneous teacher-name text with the real
# teaching method.
$metodoText = "Aulas expositivas, exercícios e visitas didátias de campo."
$ws.Range("B19").Value = $metodoText
$ws.Range("C19").Value = $metodoText

# Row 20 (Critério:): replace with the assessment criteria text.
$criterioText = "Provas e relatórios."
$ws.Range("B20").Value = $criterioText
$ws.Range("C20").Value = $criterioText

# Row 21 (Norma de recuperação:): replace with the recovery-exam rule text.
$normaText = "Prova única com nota igual ou superior a 5,0 (cinco)."
$ws.Range("B21").Value = $normaText
$ws.Range("C21").Value = $normaText

# Row 22 (Bibliografia:): set the full bibliography text.
$bibliografiaText = "ATTEWELL & FARMER - 1976 - Principles of Enginnering Geology. Chapman Hall.`nDUNICLIFF, J. - 1988 - Geotechnical Instrumentation for Monitoring Field Performance, Joh Willey & Sons, New York, 577 p.`nHANNA, T.H. - 1996 - Field Instrumentation in Geotechnical Engineering. Trans Tech Publications, RockPort - MA, 843 p.`nKELLY, W.E. e MARES S. - Applied Geophyses in Hydrogeological and Engineering Practice. Elsevier, New York - 1993, 300p.`nLUIZ, J.G. - 1995 - Geofísica de Prospecção. Editora Universitária UFPA, Belém, 1995. `nVOGELSAND, D. - 1995 - Environmental Geophysics. Springer - Verlag, Berlin, 171p."
$ws.Range("B22").Value = $bibliografiaText
$ws.Range("C22").Value = $bibliografiaText
